# Journal de travail - correction (Da Silva Diogo)
# Applies the edits described in the commit "correction journal travail diogo"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- Week 1 (rows 6-14) -------------------------------------------------
# "installation Raspberry Pi" task took 2h instead of 1h
$ws.Range("D7").Value = 2

# New task logged for the second day of the week: "doc" - 2h
$ws.Range("B11").Value = "doc"
$ws.Range("D11").Value = 2

# Weekly personal reflection text
$ws.Range("A15").Value = "Cette semaine j'ai principalement travaillé sur la compréhension du projet"

# --- Week 3 (rows 30-43) -------------------------------------------------
$ws.Range("D35").Value = 2
$ws.Range("D39").Value = 2

# --- Week 4 (rows 44-57) -------------------------------------------------
# Fix a year typo in the dates: 2025-01-05/06 -> 2026-01-05/06
$ws.Range("A44").Value = "01/05/26"
$ws.Range("A48").Value = "01/05/26"
$ws.Range("A52").Value = "01/06/26"

# Weekly personal reflection text
$ws.Range("A57").Value = "Je pense avoir réalisé un travail utile, qui permet à d'autres personnes d'installer facilement Arcadiabox là où elles le souhaitent."

# --- Recalculate & refresh the view -------------------------------------
$excel.Calculate()

$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("I50").Select()
